# Generate Report for handback
#
# Marks the two in-flight localization files as handed back (in sync with
# en-US) for both target languages (zh-cn, de-de): stamps the "Latest
# Target File" / "Latest Handback File" columns (E/F) with the same
# source/xlf file references already shown in columns A/C, refreshes the
# "Latest Handback DateTime" (column G), and flips the overall Status
# (column B, plus the Overview rollup) from "Ready for handoff" to
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: the Status rollup for both rows (zh-cn & de-de columns)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Helper data: per-language hyperlink targets, reused for both the
# existing columns (A/C) and the newly populated ones (E/F).
# ---------------------------------------------------------------------

# zh-cn sheet
$ws = $wb.Worksheets.Item("zh-cn")

$mdUrl1  = "https://github.com/OpenLocalizationTest/oltest/blob/3f139c093c6c020e9b069107aef7caec89b4bf53/e2e/3aaff79e-7311-419c-9ef3-0ea864b799da.md"
$xlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d74910d327150fa34d2b892b174ffb1eb73e82eb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.zh-cn.xlf"
$mdUrl2  = "https://github.com/OpenLocalizationTest/oltest/blob/3f139c093c6c020e9b069107aef7caec89b4bf53/e2e/9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"
$xlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d74910d327150fa34d2b892b174ffb1eb73e82eb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.zh-cn.xlf"

# Row 2 - 3aaff79e-...md source file
$ws.Range("B2").Value = $newStatus
$ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl1, "", "", "3aaff79e-7311-419c-9ef3-0ea864b799da.md")
$ws.Range("E2").Font.Underline = 2
$ws.Range("E2").Font.Color = 15570276
$ws.Hyperlinks.Add($ws.Range("F2"), $xlfUrl1, "", "", "3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.zh-cn.xlf")
$ws.Range("F2").Font.Underline = 2
$ws.Range("F2").Font.Color = 15570276
$ws.Range("G2").Value = "2016-02-19 05:45:27"
$ws.Range("H2").Value = "Include"

# Row 3 - 9eb1fb6a-...md source file
$ws.Range("B3").Value = $newStatus
$ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl2, "", "", "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md")
$ws.Range("E3").Font.Underline = 2
$ws.Range("E3").Font.Color = 15570276
$ws.Hyperlinks.Add($ws.Range("F3"), $xlfUrl2, "", "", "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.zh-cn.xlf")
$ws.Range("F3").Font.Underline = 2
$ws.Range("F3").Font.Color = 15570276
$ws.Range("G3").Value = "2016-02-19 05:45:27"
$ws.Range("H3").Value = "Include"

# de-de sheet
$ws = $wb.Worksheets.Item("de-de")

$mdUrl1  = "https://github.com/OpenLocalizationTest/oltest/blob/3f139c093c6c020e9b069107aef7caec89b4bf53/e2e/3aaff79e-7311-419c-9ef3-0ea864b799da.md"
$xlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4a43525ba630a81de43cbcf7977460be8cb2f356/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.de-de.xlf"
$mdUrl2  = "https://github.com/OpenLocalizationTest/oltest/blob/3f139c093c6c020e9b069107aef7caec89b4bf53/e2e/9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"
$xlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4a43525ba630a81de43cbcf7977460be8cb2f356/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.de-de.xlf"

# Row 2 - 3aaff79e-...md source file
$ws.Range("B2").Value = $newStatus
$ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl1, "", "", "3aaff79e-7311-419c-9ef3-0ea864b799da.md")
$ws.Range("E2").Font.Underline = 2
$ws.Range("E2").Font.Color = 15570276
$ws.Hyperlinks.Add($ws.Range("F2"), $xlfUrl1, "", "", "3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.de-de.xlf")
$ws.Range("F2").Font.Underline = 2
$ws.Range("F2").Font.Color = 15570276
$ws.Range("G2").Value = "2016-02-19 05:45:43"
$ws.Range("H2").Value = "Include"

# Row 3 - 9eb1fb6a-...md source file
$ws.Range("B3").Value = $newStatus
$ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl2, "", "", "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md")
$ws.Range("E3").Font.Underline = 2
$ws.Range("E3").Font.Color = 15570276
$ws.Hyperlinks.Add($ws.Range("F3"), $xlfUrl2, "", "", "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.de-de.xlf")
$ws.Range("F3").Font.Underline = 2
$ws.Range("F3").Font.Color = 15570276
$ws.Range("G3").Value = "2016-02-19 05:45:43"
$ws.Range("H3").Value = "Include"
